# Changelog edit: add the "0.3.0" test release entry (row 5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 4's formatting down into row 5 (same as a user dragging the
# fill handle / copying the row) so the new cells pick up the existing date,
# text and wrap-text styles instead of minting new ones.
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C5").PasteSpecial("formats")
$ws.Range("C5").WrapText = $true

# New changelog entry content
$ws.Range("A5").Value = 44293
$ws.Range("B5").Value = "0.3.0"
$ws.Range("C5").Value = "Next testversion with all content and fixed bugs"

# Leave the selection where the author ended up after typing the new row
$ws.Range("C14").Select()
